# Fix documentation about how to import the project
#
# Slide 10 ("Durchführung I"), shape "Textplatzhalter 3": the 8th
# paragraph previously read "Öffnet das Projekt in IntelliJ". It is
# replaced with a corrected description of how to import the project,
# including the (bold, Consolas) "java" source-folder name.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$para8 = $tr.Paragraphs(8, 1)
$para8.Text = "Öffne IntelliJ, importiere den java Ordner als neues Projekt"

$fullText = $para8.Text
$javaIndex = $fullText.IndexOf("java")
$javaRun = $para8.Characters($javaIndex + 1, 4)
$javaRun.Font.Bold = -1
$javaRun.Font.Name = "Consolas"
